$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.868.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -1.47%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.347.87"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -1.68%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'505.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.18%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'129.77"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -1.93%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.536"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -2.49%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.357.78"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -1.44%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -0.17%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -0.20%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'4.79"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +2.72%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.320"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.93%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'2.763.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.74%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'55.816.78"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -1.41%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E17").Value = "'  -0.56%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.317.47"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -2.30%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'9.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -2.64%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'310.44"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +0.41%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'4.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -0.98%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'6.20"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -0.83%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  -0.14%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'65.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -2.09%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.997"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.21%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  -3.59%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  -3.19%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'7.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -4.37%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'171.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -2.75%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -0.60%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.0₃0704"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -2.83%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +0.00%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'5.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -1.32%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'0.996"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -0.17%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  -4.57%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  -0.78%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  -1.87%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.830"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.69%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'3.64"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -4.14%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'36.11"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -2.00%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'1.39"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -3.63%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'3.34"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -1.08%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'4.86"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.48%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'125.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -4.13%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.556"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -1.75%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.0892"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -1.85%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'239.36"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -3.98%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.0476"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -1.30%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'16.73"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -1.50%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0206"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -1.85%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'16.69"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -2.87%  "
$ws.Range("E51").ClearFormats()
